# Daily auto push: insert a new reading for 2026/01/28 right before the
# 2026/12/29 block (row 730), pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 730; everything from the old row 730 onward
# (2026/12/29 ... 2027/01/05) shifts down to 731 onward.
$ws.Rows.Item(730).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/01/28"), not as
# real Excel dates. Temporarily force text formatting so Excel doesn't
# auto-convert the string into a date serial number, then clear the
# formatting override afterwards so the cell keeps the sheet's default
# (unstyled) look, matching the other data rows.
$ws.Range("A730").NumberFormat = "@"
$ws.Range("A730").Value = "2026/01/28"
$ws.Range("A730").ClearFormats()

$ws.Range("B730").Value = "水"
$ws.Range("C730").Value = 23
$ws.Range("D730").Value = 201
